# TC15_Verify_ViewFullCart.xlsx - "Changes done for Kaman new UI - header & footer"
#
# 1. The SCROLL_DOWN keyword step before "Add To Cart" is replaced with a
#    TINY_SCROLL_DOWN keyword (smaller scroll now suffices for the new UI).
# 2. The "CLICK ViewFullCart" step is no longer needed (new UI goes straight
#    to the cart totals), so that row is removed entirely and the rows below
#    it shift up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 14: SCROLL_DOWN -> TINY_SCROLL_DOWN
$ws.Range("B14").Value = "TINY_SCROLL_DOWN"

# Row 19 was: CLICK | ViewFullCart | CSS -- remove it, shifting rows 20:22 up
$ws.Rows.Item(19).Delete()

# Leave the selection on the (now) last data row, matching the saved file
[void]$ws.Rows.Item(19).Select()
